$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "maca"

$ws.Range("A4").Value = "maçã"
$ws.Range("B4").Value = 15

$ws.Range("A5").Value = "abacaxi"
$ws.Range("B5").Value = 20
